$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data (16th May refresh) - 3 additional reg_center_user_h records
$newRows = @(
    @(10005, 110033),
    @(10005, 110034),
    @(10005, 110035)
)

$startRow = 34
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $newRows[$i][1]
    $ws.Cells.Item($r, 3).Value = "eng"
    $ws.Cells.Item($r, 4).Value = $true
    $ws.Cells.Item($r, 5).Value = "superadmin"
    $ws.Cells.Item($r, 6).Value = "now()"
    $ws.Cells.Item($r, 7).Value = "now()"
}

# Update selection to reflect post-edit state (rows from 37 to the end of the sheet selected)
$ws.Range("A37:XFD1048576").Select()
